$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mensajes de wapp")
$ws.Activate()

# Mark the test order in row 3 as not (yet) processed: "Si" -> "No"
$ws.Range("A3").Value = "No"

# Leave the selection where it would land after editing A3 (Enter moves down one row)
$ws.Range("A4").Select()

try {
    $excel.ActiveWindow.ScrollRow = 3
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # view-only scroll position; safe to ignore if unsupported
}
